# WAT new test cases implementation
# Applies the sharedStrings / sheet1 changes described by the diff:
#  - Row 46 (WAT32 / WAT-567) description split away from the old combined
#    WAT-567||WAT-568 entry, row height reset to default.
#  - Four new rows (47-50) added for WAT33-WAT36.
#  - Selection moved from C53 to B53.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Row 46: split the old combined JIRA id / description into a single one ----
$ws.Range("B46").Value = 'WAT-567'
$ws.Range("C46").Value = 'Verify that System must display the department name (sub-organization) in addition to the institution/org name.'
$ws.Rows.Item(46).AutoFit()

# ---- New row 47 : WAT33 / WAT-611 ----
$ws.Range("A47").Value = 'WAT33'
$ws.Range("B47").Value = 'WAT-611'
$ws.Range("C47").Value = 'Verify that System must provide "Search Results" tab on Author record page'
$ws.Range("D47").Value = 'Y'

# ---- New row 48 : WAT34 / WAT-613 ----
$ws.Range("A48").Value = 'WAT34'
$ws.Range("B48").Value = 'WAT-613'
$ws.Range("C48").Value = 'Verify that User must be navigated back to the ORCID search results page and the ORCID original search results must be displayed when the "Search results" tab is clicked.'
$ws.Range("D48").Value = 'Y'
$ws.Rows.Item(48).RowHeight = 30

# ---- New row 49 : WAT35 / WAT-612 ----
$ws.Range("A49").Value = 'WAT35'
$ws.Range("B49").Value = 'WAT-612'
$ws.Range("C49").Value = 'Verify that User must be navigated back to the search results page and the original search results must be displayed when the "Search results" tab is clicked.'
$ws.Range("D49").Value = 'Y'
$ws.Rows.Item(49).RowHeight = 30

# ---- New row 50 : WAT36 / WAT-614 ----
$ws.Range("A50").Value = 'WAT36'
$ws.Range("B50").Value = 'WAT-614'
$ws.Range("C50").Value = 'Verify that "Search Results" tab should be highlighted when user navigate back from Author record page to Search Results page'
$ws.Range("D50").Value = 'Y'

# ---- Copy formatting (borders / wrap / vertical alignment) from the existing ----
# ---- styled rows so the new cells match the look of the rest of the table.   ----
$ws.Range("A46").Copy()
$ws.Range("A47:A50").PasteSpecial(-4122)

$ws.Range("B46").Copy()
$ws.Range("B47:B50").PasteSpecial(-4122)

$ws.Range("C45").Copy()
$ws.Range("C47:C50").PasteSpecial(-4122)

$ws.Range("D46").Copy()
$ws.Range("D47:D50").PasteSpecial(-4122)

$ws.Range("E46").Copy()
$ws.Range("E47:E50").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---- Update the active selection to match the authored workbook (B53) ----
$ws.Range("B53").Select()
